$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the old header value in B1 ("Jun_13") before we shift things around.
$oldB1 = $ws.Range("B1").Value2

# Insert two new columns before column C. This pushes the existing column C
# (analyst rating values/styles) two slots to the right, into column E.
$ws.Range("C1:D1").EntireColumn.Insert()

# New column headers in row 1.
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"
$ws.Range("D1").Value = $oldB1

# Fill the two newly inserted columns (rows 2-27) with the "UN" placeholder,
# matching column B's placeholder values.
$ws.Range("C2:D27").Value = "UN"

# Match column widths across C, D and E (all ~8 characters wide).
$ws.Columns("C").ColumnWidth = 7.14
$ws.Columns("D").ColumnWidth = 7.14
$ws.Columns("E").ColumnWidth = 7.14
